$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.022.14'
$ws.Range("E2").Value = '  -0.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.160.33'
$ws.Range("E3").Value = '  -0.76%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.83'
$ws.Range("E5").Value = '  +2.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.49'
$ws.Range("E6").Value = '  -2.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.160.50'
$ws.Range("E8").Value = '  -0.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  +0.30%  '

$ws.Range("E10").Value = '  -3.77%  '

$ws.Range("E11").Value = '  -0.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.503'
$ws.Range("E12").Value = '  -1.24%  '

$ws.Range("E13").Value = '  -4.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.48'
$ws.Range("E14").Value = '  -2.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.679.29'
$ws.Range("E15").Value = '  -0.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.031.04'
$ws.Range("E16").Value = '  -0.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.20'
$ws.Range("E17").Value = '  -0.76%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.162.93'
$ws.Range("E18").Value = '  -0.72%  '

$ws.Range("E19").Value = '  -0.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '505.48'
$ws.Range("E20").Value = '  -2.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.60'
$ws.Range("E21").Value = '  +4.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.721'
$ws.Range("E22").Value = '  -2.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.06'
$ws.Range("E23").Value = '  -6.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.83'
$ws.Range("E24").Value = '  -0.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.68'
$ws.Range("E25").Value = '  -0.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").Value = '  -0.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.95'
$ws.Range("E28").Value = '  +1.20%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.21'
$ws.Range("E29").Value = '  -0.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.84'
$ws.Range("E30").Value = '  +0.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.77'
$ws.Range("E31").Value = '  -1.64%  '

$ws.Range("E32").Value = '  -0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.48'
$ws.Range("E34").Value = '  +2.80%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.52'
$ws.Range("E35").Value = '  -2.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.93'
$ws.Range("E36").Value = '  -1.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0901'
$ws.Range("E37").Value = '  +2.33%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '471.78'
$ws.Range("E38").Value = '  -1.94%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0419'
$ws.Range("E39").Value = '  -1.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.93'
$ws.Range("E40").Value = '  -7.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.84'
$ws.Range("E41").Value = '  +1.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.995.58'
$ws.Range("E42").Value = '  -4.95%  '

$ws.Range("E43").Value = '  -3.57%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.285'
$ws.Range("E44").Value = '  -2.01%  '

$ws.Range("E45").Value = '  -3.64%  '

$ws.Range("B46").Value = 'PEPE'
$ws.Range("C46").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₃0610'
$ws.Range("E46").Value = '  +3.03%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.38'
$ws.Range("E47").Value = '  -3.50%  '

$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("E49").Value = '  -1.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.27'
$ws.Range("E50").Value = '  -3.00%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.32'
$ws.Range("E51").Value = '  -3.49%  '
